# Fix for GetColumns(): the TEST6 sample sheet was missing a leading blank
# row, which hid an off-by-one when the header row was detected. Insert a
# new blank row above the current row 1 so the header + data rows shift
# down by one (old row 1 -> row 2, row 2 -> row 3, row 3 -> row 4, row 4 -> row 5).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TEST6")
$ws.Activate()

$ws.Rows.Item(1).Insert()

# Leave the sheet's cell selection where the author's last save left it.
$ws.Range("M17").Select()
